$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("raw data")
$ws2 = $wb.Worksheets.Item("linear regression")

# ---------------------------------------------------------------------------
# 1. "raw data" sheet: append a new measurement row (row 9)
# ---------------------------------------------------------------------------
$ws1.Range("A9").Value = "Joshua Thomas Brooks"

# Copy the date formatting from B2 (style index 1) onto B9 so we don't create
# a brand-new (duplicate) number format in styles.xml.
$ws1.Range("B2").Copy()
$ws1.Range("B9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("B9").Value = 45509

$ws1.Range("C9").Value = 1703
$ws1.Range("D9").Value = 33
$ws1.Range("E9").Value = 1
$ws1.Range("F9").Value = 5
$ws1.Range("G9").Value = 123
$ws1.Range("H9").Value = 76
$ws1.Range("I9").Value = 96
$ws1.Range("J9").Value = 77
$ws1.Range("K9").Value = 98.2
$ws1.Range("L9").Value = 97
$ws1.Range("M9").Value = 112

# Re-create the AutoFilter so its range grows from A1:M2 to A1:M9.
$ws1.AutoFilterMode = $false
[void]$ws1.Range("A1:M9").AutoFilter()

# The hidden _FilterDatabase defined name also needs to point at the new
# range (AutoFilter() above doesn't always touch it).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "raw data!_FilterDatabase") {
        $n.RefersTo = "='raw data'!`$A`$1:`$M`$9"
    }
}

# Selection on "raw data" moves to the whole of column A.
[void]$ws1.Columns.Item(1).Select()

# ---------------------------------------------------------------------------
# 2. "linear regression" sheet: new sample inputs + refit coefficients
# ---------------------------------------------------------------------------
$ws2.Range("D2").Value = 5
$ws2.Range("E2").Value = 123
$ws2.Range("F2").Value = 76
$ws2.Range("G2").Value = 96
$ws2.Range("H2").Value = 77
$ws2.Range("I2").Value = 98.2
$ws2.Range("J2").Value = 97

$ws2.Range("B3").Value = -196383.101
$ws2.Range("B4").Value = 201.9364
$ws2.Range("B5").Value = -143.6046
$ws2.Range("B6").Value = 363.5279
$ws2.Range("B7").Value = -89.0195
$ws2.Range("B8").Value = 90.2147
$ws2.Range("B9").Value = 1535.7889
$ws2.Range("B10").Value = 374.2838

# Updated regression-formula description text.
$ws2.Range("A16").Value = "Interstitial Fluid = -196383.1010 + 201.9364 * Acetone ketones ppm - 143.6046 * Blood Pressure Top + 363.5279 * Blood Pressure Bottom - 89.0195 * Pulse + 90.2147 * Heart Rate + 1535.7889 * Temperature + 374.2838 * Pulse Oxygen"

# Selection on "linear regression" moves to C12.
[void]$ws2.Range("C12").Select()
[void]$ws2.Activate()

$wb.Save()
